$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 824
$ws1.Range("F5").Value = 844
$ws1.Range("F6").Value = 478
$ws1.Range("F7").Value = 630
$ws1.Range("F8").Value = 195
$ws1.Range("F9").Value = 8
$ws1.Range("F12").Value = 118
$ws1.Range("G12").Value = 29.9
$ws1.Range("F13").Value = 1626
$ws1.Range("F14").Value = 206
$ws1.Range("F17").Value = 74
$ws1.Range("F18").Value = 395
$ws1.Range("F23").Value = 732
$ws1.Range("F25").Value = 1463
$ws1.Range("F26").Value = 184

# Sheet "演出" updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 272

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 824
$ws4.Range("F6").Value = 844
$ws4.Range("F9").Value = 478
$ws4.Range("F10").Value = 630
$ws4.Range("F12").Value = 195
$ws4.Range("F13").Value = 8
$ws4.Range("F16").Value = 118
$ws4.Range("G16").Value = 29.9
$ws4.Range("F17").Value = 1626
$ws4.Range("F19").Value = 206
$ws4.Range("F22").Value = 74
$ws4.Range("F23").Value = 395
$ws4.Range("F26").Value = 272
$ws4.Range("F35").Value = 732
$ws4.Range("F37").Value = 1463
$ws4.Range("F38").Value = 184
